$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6784462332725525
$ws.Range("B1").Value = 2.512659311294556
$ws.Range("C1").Value = 3.282250642776489
$ws.Range("D1").Value = 3.699621200561523
$ws.Range("E1").Value = 0.9225375652313232
